# CDS Input file updates
#
# The ParticipantsTab query (cell B2 on the "startup" sheet) is replaced with
# a new Cypher query that also optionally matches diagnosis/genomic_info and
# returns the sample ids sorted via apoc.coll.sort. The other query cells
# (C2/B3/C3/B4/C4) are untouched - Excel's own shared-string bookkeeping
# accounts for the index renumbering seen in the saved file.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newQuery = @"
MATCH (p:participant)-->(s:study)
OPTIONAL MATCH (samp:sample)-->(p)
OPTIONAL MATCH (p)<--(diag:diagnosis)
OPTIONAL MATCH (samp)<--(f:file)
OPTIONAL MATCH (f)<--(g:genomic_info)
WITH s, p, samp, f, g, diag
WHERE f.file_type in ['IDAT']
with p
OPTIONAL MATCH (p)-->(s:study)
OPTIONAL MATCH (samp:sample)-->(p)
WITH s, p, apoc.coll.sort(collect(distinct samp.sample_id)) as samp
RETURN
coalesce(p.participant_id,'') as ``Participant ID``,
coalesce(s.study_name, '') as ``Study Name``,
coalesce(s.phs_accession,'') as ``Accession``,
coalesce(p.gender,'') as ``Gender``,
coalesce(apoc.text.join(samp, ','), '') as ``Samples``
ORDER BY p.participant_id LIMIT 100
"@

$ws.Range("B2").Value = $newQuery

# The new query text wraps across more lines than the old one, so row 2 grows
# from 186 to 279 points.
$ws.Rows.Item(2).RowHeight = 279

# The saved view had the selection moved to B5 (with the sheet scrolled so
# row 3 is at the top).
$ws.Range("B5").Select() | Out-Null
